# [Kadastro App] Yeni kayit eklendi: 1
#
# Appends one new record row (row 2) to the "Kayitlar" (master log) sheet
# and to the "Erdemli" (per-birim) sheet - both currently hold only the
# header row. The new row's values must land as literal TEXT (matching the
# existing header cells and the workbook's numberStoredAsText convention),
# even though some of them look numeric/date-like, so the target range is
# pre-formatted as Text before the values are typed in.

$wb = $excel.ActiveWorkbook

$kayitNo    = "1"
$tarih      = "2025-09-04"
$birim      = "Erdemli"
$parselSayi = "50"
$is         = "18-UYG."
$personel   = "EMİNE ALANLI KIRCILI (K.Mühendisi), MAHMUT GÖK (Kontrol Memuru)"

$rowValues = @($kayitNo, $tarih, $birim, $parselSayi, $is, $personel)

foreach ($sheetName in @("Kayitlar", "Erdemli")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $targetRange = $ws.Range("A2:F2")
    # Force text storage so numeric-looking values ("1", "2025-09-04", "50")
    # are kept as strings instead of being auto-converted to number/date.
    $targetRange.NumberFormat = "@"

    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item(2, $col).Value = $rowValues[$col - 1]
    }
}
